$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F3").Value = 223
$ws.Range("F5").Value = 9009
$ws.Range("F8").Value = 146
$ws.Range("F9").Value = 181
$ws.Range("F10").Value = 320
$ws.Range("F15").Value = 11496
$ws.Range("F16").Value = 11496
$ws.Range("F34").Value = 48
$ws.Range("F37").Value = 944
$ws.Range("F38").Value = 4138
$ws.Range("F39").Value = 306
$ws.Range("F41").Value = 1277
$ws.Range("F44").Value = 372
$ws.Range("F45").Value = 406

$ws = $wb.Worksheets.Item(2)
$ws.Range("F7").Value = 15
$ws.Range("F9").Value = 37
$ws.Range("F14").Value = 17
$ws.Range("G20").Value = 280

$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 14

$ws = $wb.Worksheets.Item(4)
$ws.Range("F7").Value = 223
$ws.Range("F9").Value = 9009
$ws.Range("F11").Value = 37
$ws.Range("F13").Value = 181
$ws.Range("F14").Value = 320
$ws.Range("F18").Value = 11496
$ws.Range("F35").Value = 48
$ws.Range("F38").Value = 944
$ws.Range("F40").Value = 4138
$ws.Range("F41").Value = 306
$ws.Range("F43").Value = 1277
$ws.Range("F45").Value = 372
